try {
  $t = [System.Reflection.Assembly]
  Write-Host "reflection available: $t"
} catch {
  Write-Host "fail: $_"
}
